$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows at position 9 (pushes the old totals row 9 -> 11
#    and the old footer row 10 -> 12), restricted to columns A:Q so we don't
#    touch/format the whole 16384-column row.
# ---------------------------------------------------------------------------
$ws.Range("A9:Q10").Insert()

# ---------------------------------------------------------------------------
# 2. Copy the cell formatting of the two existing item rows (7 and 8) onto
#    the two freshly inserted rows (9 and 10) so they look like item rows.
# ---------------------------------------------------------------------------
$ws.Range("A7:Q8").Copy()
$ws.Range("A9:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights for the new rows (matching the other item rows).
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75

# Re-create the merged cells for the new item rows.
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# ---------------------------------------------------------------------------
# 3. New item #3 : CATAFAST 50 MG 9 SACHET
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "CATAFAST 50 MG 9 SACHET"
$ws.Range("H9").Value = "1:1"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "72.00"
$ws.Range("P9").Value = "7.9200"
$ws.Range("Q9").Value = "0:1"

# ---------------------------------------------------------------------------
# 4. New item #4 : TORSERETIC 20MG 30 TABS. (was previously item #2)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "TORSERETIC 20MG 30 TABS."
$ws.Range("H10").Value = "0:2"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "123.00"
$ws.Range("P10").Value = "123.0000"
$ws.Range("Q10").Value = "1:0"

# ---------------------------------------------------------------------------
# 5. Item #2 (row 8) becomes CALMEPAM 3MG 20 TAB.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "CALMEPAM 3MG 20 TAB."
$ws.Range("H8").Value = "0:0"
$ws.Range("N8").Value = "76.00"
$ws.Range("P8").Value = "76.0000"
$ws.Range("Q8").Value = "1:0"

# ---------------------------------------------------------------------------
# 6. Totals row (now row 11) - update the grand total of the "sell price"
#    column to reflect the two new items.
# ---------------------------------------------------------------------------
$ws.Range("P11").Value = 226.72

# ---------------------------------------------------------------------------
# 7. Footer row (now row 12) - refresh the generated timestamp.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Wednesday, 23 July, 2025 10:11 AM"
